# Alteracao - Telas novas 1
# Remove the "EXCLUIR" procedure rows for FERIAS, ACIDENTE_TRABALHO and
# ALTERACAO_CARGO_SALARIO, then renumber the remaining SPVRT sequence
# numbers (both the SEQ column and the embedded number in NOME PROC) so
# the list stays contiguous.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the three rows (bottom-to-top so row numbers of the
# not-yet-deleted rows stay valid).
$ws.Rows(76).Delete()   # SPVRT075_ALTERACAO_CARGO_SALARIO_PR_EXCLUIR
$ws.Rows(72).Delete()   # SPVRT071_ACIDENTE_TRABALHO_PR_EXCLUIR
$ws.Rows(68).Delete()   # SPVRT067_FERIAS_PR_EXCLUIR

# After the deletions, rows 68..75 hold what used to be
# SPVRT068, 069, 070, 072, 073, 074, 076, 077 - renumber them
# sequentially as SPVRT067..SPVRT074 and fix the SEQ column (A).
$names = @(
  "SPVRT067_FERIAS_PR_INCLUIR",
  "SPVRT068_FERIAS_PR_SELECIONAR",
  "SPVRT069_ACIDENTE_TRABALHO_PR_ALTERAR",
  "SPVRT070_ACIDENTE_TRABALHO_PR_INCLUIR",
  "SPVRT071_ACIDENTE_TRABALHO_PR_SELECIONAR",
  "SPVRT072_ALTERACAO_CARGO_SALARIO_PR_ALTERAR",
  "SPVRT073_ALTERACAO_CARGO_SALARIO_PR_INCLUIR",
  "SPVRT074_ALTERACAO_CARGO_SALARIO_PR_SELECIONAR"
)

$startRow = 68
$startSeq = 67
for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $startRow + $i
    $seq = $startSeq + $i
    $ws.Range("A$row").Value = $seq
    $ws.Range("B$row").Value = $names[$i]
}

# Update the view/selection left over from editing near the bottom of the
# (now shorter) table.
$ws.Range("B76").Select()
[void]($excel.ActiveWindow.ScrollRow = 61)
[void]($excel.ActiveWindow.ScrollColumn = 1)

